# Add the new "stheVS26" kinetic model entry to row 17 (the
# Streptococcus thermophilus / STH_CIRM_65 placeholder row) and bring the
# whole sheet's cell alignment in line with the rest of the table
# (vertical-center), mirroring the author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Models")

# --- Fill in the new model row (row 17) ---------------------------------
$ws.Range("A17").Value = "stheVS26"
$ws.Range("B17").Value = 2026
$ws.Range("E17").Value = 336
$ws.Range("F17").Value = 400
$ws.Range("G17").Value = 322

# --- Normalize vertical alignment across the table -----------------------
# (center vertically, matching the rest of the populated cells; the engine's
# multi-area Range setter only honors the first area, so apply cell by cell)
$vCenterCells = @(
    "A13","B13","D13","F13","G13",
    "A14","B14","D14","F14","G14",
    "A15","B15","E15","F15","G15",
    "D16",
    "A17","B17","D17","E17","F17","G17",
    "A18","B18","D18","E18","F18","G18",
    "D19",
    "A20","B20","D20","E20","F20","G20",
    "D21",
    "D22",
    "A23","B23","D23","E23","F23","G23",
    "D24","D25","D26","D27","D28"
)
foreach ($addr in $vCenterCells) {
    $ws.Range($addr).VerticalAlignment = -4108
}

$italicCells = @("C16","C17","C18","C19","C20","C21","C22","C23","C24","C25","C26","C27","C28")
foreach ($addr in $italicCells) {
    $ws.Range($addr).VerticalAlignment = -4108
}

$ws.Range("H2:H15").VerticalAlignment = -4108

# --- Reposition the active cell/selection like the saved workbook --------
$ws.Range("D14").Select() | Out-Null
